$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2022", 2
)

# 2. Fill in the empty italic "Objetivos" (English) paragraph with its text
$pObjetivosEn = $d.Paragraphs.Item(7)
$pObjetivosEn.Range.Text = "Provide students with the knowledge of cell biology necessary to understand the other subjects of the course and the training of the Environmental Engineer."

# 3. "Programa resumido" - Portuguese: drop the leading "Origem e evolução da célula; " clause
$pResumoPt = $d.Paragraphs.Item(11)
$pResumoPt.Range.Text = "Análise estrutural das células ao microscópio; moléculas orgânicas; organização interna da célula; organelas celulares transdutoras de energia; material genético e mecanismo de divisão celular."

# 4. "Programa resumido" - English: drop the leading "The origin and evolution of the cell; " clause
$pResumoEn = $d.Paragraphs.Item(12)
$pResumoEn.Range.Text = "Organic molecules; internal organization of the cell; cell energy conversion; genetic material and mechanism of cell division."

# 5. "Programa" - Portuguese: replace the first bullet item
$pProgramaPt = $d.Paragraphs.Item(14)
$pProgramaPt.Range.Text = "- Estrutura celular e história evolutiva: microrganismos procarióticos eeucarióticos e suas relações evolutivas dentre os domínios Bacteria, Archaea eEukarya.- Análise estrutural das células ao microscópio: microscopia ótica e eletrônica.- Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos, ácidos nucleicos, aminoácidos. - Organização interna da célula: estrutura e função da membrana plasmática; compartimentos intracelulares e seleção de proteínas; tráfico de vesículas (via de exocitose e endocitose).- Núcleo e organização do material genético: estrutura e função- Ciclo celular e divisão celular: mitose e meiose.- Organelas celulares transdutoras de energia: mitocôndria e cloroplasto."

# 6. "Programa" - English: replace the first bullet item
$pProgramaEn = $d.Paragraphs.Item(15)
$pProgramaEn.Range.Text = "Cell structure and evolutionary history: prokaryotic microorganisms andeukaryotic and their evolutionary relationships between the Bacteria, Archaea andEukarya.Microscope analysis of cells structure: optical and electron microscope.Structure and function of major organic molecules: carbohydrates, lipids, nucleic acids and proteins. Internal organization of the cell: membrane structure and function; intracelular compartments and protein sorting; vesicular traffic (endocytosis and exocytosis).Nucleus and genetic material organization: structure and functionCell cycle and cell division: mitosis and meiosisCell energy conversion: mitochondria and chloroplast."
